$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old summary formula rows 15-17 (they move to 29-31)
$ws.Range("K15:L17").ClearContents()

# Force column F (Win %) to Text format so literal "NN%" strings are not
# auto-converted to percentage numbers by the COM layer
$ws.Range("F2:F28").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Sat Oct 11"
$ws.Range("B2").Value = "Norway ✓ - Israel: 5:0"
$ws.Range("C2").Value = 4.79
$ws.Range("D2").Value = "Norway"
$ws.Range("E2").Value = 5.5
$ws.Range("F2").Value = "80%"
$ws.Range("G2").Value = "✓"
$ws.Range("H2").Value = 5
$ws.Range("I2").Value = $true

# Row 3
$ws.Range("A3").Value = "Sat Oct 11"
$ws.Range("B3").Value = "Portugal ✓ - Republic of Ireland: 1:0"
$ws.Range("C3").Value = 3.02
$ws.Range("D3").Value = "Portugal"
$ws.Range("E3").Value = 4.5
$ws.Range("F3").Value = "77%"
$ws.Range("G3").Value = "✓"
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = $true

# Row 4
$ws.Range("A4").Value = "Sat Oct 11"
$ws.Range("B4").Value = "Spain ✓ - Georgia: 2:0"
$ws.Range("C4").Value = 4.15
$ws.Range("D4").Value = "Spain"
$ws.Range("E4").Value = 5.5
$ws.Range("F4").Value = "74%"
$ws.Range("G4").Value = "✓"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = $true

# Row 5
$ws.Range("A5").Value = "Sat Oct 11"
$ws.Range("B5").Value = "Argentina ✓ - Venezuela: 1:0"
$ws.Range("C5").Value = 2.3
$ws.Range("D5").Value = "Argentina"
$ws.Range("E5").Value = 3.5
$ws.Range("F5").Value = "74%"
$ws.Range("G5").Value = "✓"
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = $true

# Row 6
$ws.Range("A6").Value = "Sat Oct 11"
$ws.Range("B6").Value = "Sociedade Esportiva Palmeiras ✓ - Esporte Clube Juventude: 4:1"
$ws.Range("C6").Value = 2.52
$ws.Range("D6").Value = "Sociedade Esportiva Palmeiras"
$ws.Range("E6").Value = 3.5
$ws.Range("F6").Value = "73%"
$ws.Range("G6").Value = "✓"
$ws.Range("H6").Value = 5
$ws.Range("I6").Value = $false

# Row 7
$ws.Range("A7").Value = "Sat Oct 11"
$ws.Range("B7").Value = "The New Saints ✓ - Penybont FC: 6:2"
$ws.Range("C7").Value = 2.37
$ws.Range("D7").Value = "The New Saints"
$ws.Range("E7").Value = 3.5
$ws.Range("F7").Value = "71%"
$ws.Range("G7").Value = "✓"
$ws.Range("H7").Value = 8
$ws.Range("I7").Value = $false

# Row 8
$ws.Range("A8").Value = "Sat Oct 11"
$ws.Range("B8").Value = "Jaguares de Córdoba ✓ - Leones FC: 4:0"
$ws.Range("C8").Value = 1.39
$ws.Range("D8").Value = "Jaguares de Córdoba"
$ws.Range("E8").Value = 2.5
$ws.Range("F8").Value = "71%"
$ws.Range("G8").Value = "✓"
$ws.Range("H8").Value = 4
$ws.Range("I8").Value = $false

# Row 9
$ws.Range("A9").Value = "Sat Oct 11"
$ws.Range("B9").Value = "Makedonikos Neapolis - POT Iraklis Thessaloniki ✓: 1:4"
$ws.Range("C9").Value = 2.47
$ws.Range("D9").Value = "POT Iraklis Thessaloniki"
$ws.Range("E9").Value = 3.5
$ws.Range("F9").Value = "71%"
$ws.Range("G9").Value = "✓"
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = $false

# Row 10
$ws.Range("A10").Value = "Sat Oct 11"
$ws.Range("B10").Value = "Bars Karakol ✓ - Asiagoal Bishkek: 6:2"
$ws.Range("C10").Value = 1.87
$ws.Range("D10").Value = "Bars Karakol"
$ws.Range("E10").Value = 2.5
$ws.Range("F10").Value = "70%"
$ws.Range("G10").Value = "✓"
$ws.Range("H10").Value = 8
$ws.Range("I10").Value = $false

# Row 11
$ws.Range("A11").Value = "Sat Oct 11"
$ws.Range("B11").Value = "Canada X - Australia: 0:1"
$ws.Range("C11").Value = 2.12
$ws.Range("D11").Value = "Canada"
$ws.Range("E11").Value = 3.5
$ws.Range("F11").Value = "69%"
$ws.Range("G11").Value = "X"
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = $true

# Row 12
$ws.Range("A12").Value = "Sat Oct 11"
$ws.Range("B12").Value = "PS Kalamata  - AO Egaleo: 1:1"
$ws.Range("C12").Value = 1.95
$ws.Range("D12").Value = "PS Kalamata"
$ws.Range("E12").Value = 2.5
$ws.Range("F12").Value = "69%"
$ws.Range("G12").ClearContents()
$ws.Range("H12").Value = 2
$ws.Range("I12").Value = $true

# Row 13
$ws.Range("A13").Value = "Sat Oct 11"
$ws.Range("B13").Value = "United Arab Emirates ✓ - Oman: 2:1"
$ws.Range("C13").Value = 1.0
$ws.Range("D13").Value = "United Arab Emirates"
$ws.Range("E13").Value = 1.5
$ws.Range("F13").Value = "67%"
$ws.Range("G13").Value = "✓"
$ws.Range("H13").Value = 3
$ws.Range("I13").Value = $false

# Row 14
$ws.Range("A14").Value = "Sat Oct 11"
$ws.Range("B14").Value = "Bulgaria - Turkiye ✓: 1:6"
$ws.Range("C14").Value = 2.05
$ws.Range("D14").Value = "Turkiye"
$ws.Range("E14").Value = 3.5
$ws.Range("F14").Value = "64%"
$ws.Range("G14").Value = "✓"
$ws.Range("H14").Value = 7
$ws.Range("I14").Value = $false

# Row 15
$ws.Range("A15").Value = "Sat Oct 11"
$ws.Range("B15").Value = "Alga Bishkek - FK Dordoi Bishkek : 14:00"
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = "FK Dordoi Bishkek"
$ws.Range("E15").Value = 0.5
$ws.Range("F15").Value = "64%"
$ws.Range("G15").ClearContents()
$ws.Range("H15").Value = 14
$ws.Range("I15").Value = $false

# Row 16
$ws.Range("A16").Value = "Sat Oct 11"
$ws.Range("B16").Value = "NK Krsko - NK Triglav Kranj ✓: 0:4"
$ws.Range("C16").Value = 4.09
$ws.Range("D16").Value = "NK Triglav Kranj"
$ws.Range("E16").Value = 5.5
$ws.Range("F16").Value = "63%"
$ws.Range("G16").Value = "✓"
$ws.Range("H16").Value = 4
$ws.Range("I16").Value = $true

# Row 17
$ws.Range("A17").Value = "Sat Oct 11"
$ws.Range("B17").Value = "Muras United Dzhalal-Abad  - Talant Besh-Küngöy: 14:00"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = "Muras United Dzhalal-Abad"
$ws.Range("E17").Value = 0.5
$ws.Range("F17").Value = "62%"
$ws.Range("G17").ClearContents()
$ws.Range("H17").Value = 14
$ws.Range("I17").Value = $false

# Row 18
$ws.Range("A18").Value = "Sat Oct 11"
$ws.Range("B18").Value = "Estonia - Italy ✓: 1:3"
$ws.Range("C18").Value = 2.74
$ws.Range("D18").Value = "Italy"
$ws.Range("E18").Value = 3.5
$ws.Range("F18").Value = "61%"
$ws.Range("G18").Value = "✓"
$ws.Range("H18").Value = 4
$ws.Range("I18").Value = $false

# Row 19
$ws.Range("A19").Value = "Sat Oct 11"
$ws.Range("B19").Value = "CF Esperança d'Andorra - Inter Club d'Escaldes ✓: 1:2"
$ws.Range("C19").Value = 5.75
$ws.Range("D19").Value = "Inter Club d'Escaldes"
$ws.Range("E19").Value = 6.5
$ws.Range("F19").Value = "60%"
$ws.Range("G19").Value = "✓"
$ws.Range("H19").Value = 3
$ws.Range("I19").Value = $true

# Row 20
$ws.Range("A20").Value = "Sat Oct 11"
$ws.Range("B20").Value = "NK Jesenice - NK Brinje Grosuplje ✓: 2:4"
$ws.Range("C20").Value = 3.27
$ws.Range("D20").Value = "NK Brinje Grosuplje"
$ws.Range("E20").Value = 4.5
$ws.Range("F20").Value = "60%"
$ws.Range("G20").Value = "✓"
$ws.Range("H20").Value = 6
$ws.Range("I20").Value = $false

# Row 21
$ws.Range("A21").Value = "Sat Oct 11"
$ws.Range("B21").Value = "FC Orsha - FC Baranovichi ✓: 0:2"
$ws.Range("C21").Value = 3.5
$ws.Range("D21").Value = "FC Baranovichi"
$ws.Range("E21").Value = 4.5
$ws.Range("F21").Value = "59%"
$ws.Range("G21").Value = "✓"
$ws.Range("H21").Value = 2
$ws.Range("I21").Value = $true

# Row 22
$ws.Range("A22").Value = "Sat Oct 11"
$ws.Range("B22").Value = "Guangdong GZ-Power ✓ - Nanjing City: 4:1"
$ws.Range("C22").Value = 2.37
$ws.Range("D22").Value = "Guangdong GZ-Power"
$ws.Range("E22").Value = 3.5
$ws.Range("F22").Value = "58%"
$ws.Range("G22").Value = "✓"
$ws.Range("H22").Value = 5
$ws.Range("I22").Value = $false

# Row 23
$ws.Range("A23").Value = "Sat Oct 11"
$ws.Range("B23").Value = "Flint Mountain - Airbus UK Broughton ✓: 0:2"
$ws.Range("C23").Value = 3.11
$ws.Range("D23").Value = "Airbus UK Broughton"
$ws.Range("E23").Value = 4.5
$ws.Range("F23").Value = "58%"
$ws.Range("G23").Value = "✓"
$ws.Range("H23").Value = 2
$ws.Range("I23").Value = $true

# Row 24
$ws.Range("A24").Value = "Sat Oct 11"
$ws.Range("B24").Value = "UD Leiria  - UD Oliveirense: 0:0"
$ws.Range("C24").Value = 1.97
$ws.Range("D24").Value = "UD Leiria"
$ws.Range("E24").Value = 2.5
$ws.Range("F24").Value = "57%"
$ws.Range("G24").ClearContents()
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = $true

# Row 25
$ws.Range("A25").Value = "Sat Oct 11"
$ws.Range("B25").Value = "Serbia X - Albania: 0:1"
$ws.Range("C25").Value = 0
$ws.Range("D25").Value = "Serbia"
$ws.Range("E25").Value = 0.5
$ws.Range("F25").Value = "57%"
$ws.Range("G25").Value = "X"
$ws.Range("H25").Value = 1
$ws.Range("I25").Value = $false

# Row 26
$ws.Range("A26").Value = "Sat Oct 11"
$ws.Range("B26").Value = "Llandudno FC ✓ - Caersws FC: 2:0"
$ws.Range("C26").Value = 2.99
$ws.Range("D26").Value = "Llandudno FC"
$ws.Range("E26").Value = 3.5
$ws.Range("F26").Value = "56%"
$ws.Range("G26").Value = "✓"
$ws.Range("H26").Value = 2
$ws.Range("I26").Value = $true

# Row 27
$ws.Range("A27").Value = "Sat Oct 11"
$ws.Range("B27").Value = "SC Cambuur Leeuwarden ✓ - De Graafschap Doetinchem: 2:0"
$ws.Range("C27").Value = 1.71
$ws.Range("D27").Value = "SC Cambuur Leeuwarden"
$ws.Range("E27").Value = 2.5
$ws.Range("F27").Value = "55%"
$ws.Range("G27").Value = "✓"
$ws.Range("H27").Value = 2
$ws.Range("I27").Value = $true

# Row 28
$ws.Range("A28").Value = "Sat Oct 11"
$ws.Range("B28").Value = "FK Laktasi ✓ - FK Slavija Sarajevo: 2:0"
$ws.Range("C28").Value = 2.25
$ws.Range("D28").Value = "FK Laktasi"
$ws.Range("E28").Value = 3.5
$ws.Range("F28").Value = "55%"
$ws.Range("G28").Value = "✓"
$ws.Range("H28").Value = 2
$ws.Range("I28").Value = $true

# Summary formulas moved to rows 29-31
$ws.Range("K29").Formula = "=COUNTIF(I:I,TRUE)"
$ws.Range("L29").Formula = "=(K29/K31)*100"
$ws.Range("K30").Formula = "=COUNTIF(I:I,FALSE)"
$ws.Range("K31").Formula = "=K29+K30"
